$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Marking row: right-answer mark value (B11) 3 -> 5
$ws.Range("B11").Value = 5

# Update Total row: total marks obtained (B12) 42 -> 70
$ws.Range("B12").Value = 70

# Update Total row: corrected/total marks text (E12) "41/84" -> "70/140"
$ws.Range("E12").Value = "70/140"
